$d = $word.ActiveDocument
$rng = $d.Content
$rng.Start = 0
$found = $rng.Find.Execute("92÷5=18, 2", $true, $false, $false, $false, $false, $true, 1, $false, "72÷4=18, 0", 1)
if (-not $found) { Write-Host "FAILED to find occurrence #0: 92÷5=18, 2" }
$found = $rng.Find.Execute("25÷4=6, 1", $true, $false, $false, $false, $false, $true, 1, $false, "75÷5=15, 0", 1)
if (-not $found) { Write-Host "FAILED to find occurrence #1: 25÷4=6, 1" }
$found = $rng.Find.Execute("38÷9=4, 2", $true, $false, $false, $false, $false, $true, 1, $false, "84÷7=12, 0", 1)
if (-not $found) { Write-Host "FAILED to find occurrence #2: 38÷9=4, 2" }
$found = $rng.Find.Execute("41÷6=6, 5", $true, $false, $false, $false, $false, $true, 1, $false, "54÷9=6, 0", 1)
if (-not $found) { Write-Host "FAILED to find occurrence #3: 41÷6=6, 5" }
$found = $rng.Find.Execute("84÷8=10, 4", $true, $false, $false, $false, $false, $true, 1, $false, "42÷2=21, 0", 1)
if (-not $found) { Write-Host "FAILED to find occurrence #4: 84÷8=10, 4" }
$found = $rng.Find.Execute("92÷3=30, 2", $true, $false, $false, $false, $false, $true, 1, $false, "48÷5=9, 3", 1)
if (-not $found) { Write-Host "FAILED to find occurrence #5: 92÷3=30, 2" }
$found = $rng.Find.Execute("64÷9=7, 1", $true, $false, $false, $false, $false, $true, 1, $false, "94÷4=23, 2", 1)
if (-not $found) { Write-Host "FAILED to find occurrence #6: 64÷9=7, 1" }
$found = $rng.Find.Execute("44÷3=14, 2", $true, $false, $false, $false, $false, $true, 1, $false, "62÷6=10, 2", 1)
if (-not $found) { Write-Host "FAILED to find occurrence #7: 44÷3=14, 2" }
$found = $rng.Find.Execute("24÷6=4, 0", $true, $false, $false, $false, $false, $true, 1, $false, "62÷8=7, 6", 1)
if (-not $found) { Write-Host "FAILED to find occurrence #8: 24÷6=4, 0" }
$found = $rng.Find.Execute("32÷2=16, 0", $true, $false, $false, $false, $false, $true, 1, $false, "48÷2=24, 0", 1)
if (-not $found) { Write-Host "FAILED to find occurrence #9: 32÷2=16, 0" }
$found = $rng.Find.Execute("35÷6=5, 5", $true, $false, $false, $false, $false, $true, 1, $false, "80÷8=10, 0", 1)
if (-not $found) { Write-Host "FAILED to find occurrence #10: 35÷6=5, 5" }
$found = $rng.Find.Execute("39÷5=7, 4", $true, $false, $false, $false, $false, $true, 1, $false, "29÷3=9, 2", 1)
if (-not $found) { Write-Host "FAILED to find occurrence #11: 39÷5=7, 4" }
$found = $rng.Find.Execute("45÷2=22, 1", $true, $false, $false, $false, $false, $true, 1, $false, "32÷5=6, 2", 1)
if (-not $found) { Write-Host "FAILED to find occurrence #12: 45÷2=22, 1" }
$found = $rng.Find.Execute("62÷8=7, 6", $true, $false, $false, $false, $false, $true, 1, $false, "93÷4=23, 1", 1)
if (-not $found) { Write-Host "FAILED to find occurrence #13: 62÷8=7, 6" }
$found = $rng.Find.Execute("10÷9=1, 1", $true, $false, $false, $false, $false, $true, 1, $false, "88÷3=29, 1", 1)
if (-not $found) { Write-Host "FAILED to find occurrence #14: 10÷9=1, 1" }
$found = $rng.Find.Execute("46÷3=15, 1", $true, $false, $false, $false, $false, $true, 1, $false, "63÷4=15, 3", 1)
if (-not $found) { Write-Host "FAILED to find occurrence #15: 46÷3=15, 1" }
$found = $rng.Find.Execute("95÷8=11, 7", $true, $false, $false, $false, $false, $true, 1, $false, "66÷2=33, 0", 1)
if (-not $found) { Write-Host "FAILED to find occurrence #16: 95÷8=11, 7" }
$found = $rng.Find.Execute("78÷5=15, 3", $true, $false, $false, $false, $false, $true, 1, $false, "20÷5=4, 0", 1)
if (-not $found) { Write-Host "FAILED to find occurrence #17: 78÷5=15, 3" }
$found = $rng.Find.Execute("38÷9=4, 2", $true, $false, $false, $false, $false, $true, 1, $false, "29÷3=9, 2", 1)
if (-not $found) { Write-Host "FAILED to find occurrence #18: 38÷9=4, 2" }
$found = $rng.Find.Execute("17÷7=2, 3", $true, $false, $false, $false, $false, $true, 1, $false, "32÷4=8, 0", 1)
if (-not $found) { Write-Host "FAILED to find occurrence #19: 17÷7=2, 3" }
$found = $rng.Find.Execute("37÷8=4, 5", $true, $false, $false, $false, $false, $true, 1, $false, "33÷9=3, 6", 1)
if (-not $found) { Write-Host "FAILED to find occurrence #20: 37÷8=4, 5" }
$found = $rng.Find.Execute("24÷6=4, 0", $true, $false, $false, $false, $false, $true, 1, $false, "41÷4=10, 1", 1)
if (-not $found) { Write-Host "FAILED to find occurrence #21: 24÷6=4, 0" }
$found = $rng.Find.Execute("51÷9=5, 6", $true, $false, $false, $false, $false, $true, 1, $false, "17÷5=3, 2", 1)
if (-not $found) { Write-Host "FAILED to find occurrence #22: 51÷9=5, 6" }
$found = $rng.Find.Execute("68÷6=11, 2", $true, $false, $false, $false, $false, $true, 1, $false, "22÷6=3, 4", 1)
if (-not $found) { Write-Host "FAILED to find occurrence #23: 68÷6=11, 2" }
$found = $rng.Find.Execute("55÷6=9, 1", $true, $false, $false, $false, $false, $true, 1, $false, "66÷3=22, 0", 1)
if (-not $found) { Write-Host "FAILED to find occurrence #24: 55÷6=9, 1" }
